# "Estadisticos Segundo Parcial 23 Mayo"
#
# 1) Fix the "Estadisticos 2P" (Second Partial) summary row for the
#    6APV group: the pass/fail counts were swapped and the percentage /
#    average columns were missing, so recompute them.
# 2) Add the two "Rescatable" (make-up exam) students that were found
#    for the 2APV group to the "Rescatables" sheet.

$wb = $excel.ActiveWorkbook

# --- 1) Estadisticos 2P ----------------------------------------------
$ws2 = $wb.Worksheets.Item("Estadisticos 2P")

$ws2.Range("D3").Value = 0      # Blancos
$ws2.Range("E3").Value = 0      # Reprobados
$ws2.Range("F3").Value = 17     # Aprobados
$ws2.Range("G3").Value = 100    # Por_Apro
$ws2.Range("H3").Value = 9.5    # Promedio

# --- 2) Rescatables ----------------------------------------------------
$ws4 = $wb.Worksheets.Item("Rescatables")

# Fill column by column so new names are registered in the same order
# the original workbook used (Paterno, Materno, Nombres, ...).
$ws4.Range("A2").Value = 24330051920274
$ws4.Range("A3").Value = 24330051920369

$ws4.Range("B2").Value = "CLEMENTE"
$ws4.Range("B3").Value = "TORRES"

$ws4.Range("C2").Value = "JUAREZ"
$ws4.Range("C3").Value = "GUTIERREZ"

$ws4.Range("D2").Value = "BRYAN"
$ws4.Range("D3").Value = "JESUS ENRIQUE"

$ws4.Range("E2").Value = "DISEÑA SOFTWARE DE SISTEMAS INFORMÁTICOS"
$ws4.Range("E3").Value = "DISEÑA SOFTWARE DE SISTEMAS INFORMÁTICOS"

$ws4.Range("F2").Value = "2APV"
$ws4.Range("F3").Value = "2APV"

$ws4.Range("G2").Value = 4
$ws4.Range("G3").Value = 3
